# "Colocando header nos gráficos"
# Adds a header label in column A (row 1) for each data table used by the
# charts, fixes accented Portuguese text that was previously stored without
# diacritics, clears the old bold/border header style off the now-plain data
# rows in column A, removes the unused "Teto" row from the Emissoes sheet,
# and updates the Custo Total sheet (new header row + new cost figures).

$wb = $excel.ActiveWorkbook

# Helper: copy the formatting of a known header cell (B1, already styled
# bold/centered/bordered) onto the newly added A1 header cell so it reuses
# the exact same style definition instead of minting a near-duplicate one.
function Set-HeaderCell($ws, $text) {
    $ws.Range("A1").Value = $text
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# ---------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
# (MWMed)", "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# All four share the same row layout (Fonte/Tecnologia table).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    Set-HeaderCell $ws "Fonte/Tecnologia"

    $ws.Range("A2").ClearFormats()

    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A3").ClearFormats()

    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A4").ClearFormats()

    $ws.Range("A5").ClearFormats()

    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A6").ClearFormats()

    $ws.Range("A7").ClearFormats()

    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A8").ClearFormats()

    $ws.Range("A9").ClearFormats()

    $ws.Range("A10").ClearFormats()

    $ws.Range("A11").Value = "Pot. Compl."
    $ws.Range("A11").ClearFormats()

    $ws.Range("A12").ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").ClearFormats()

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").ClearFormats()

# Remove the unused "Teto" row entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "Tipo Expansão"

# B1 needs to become the *text* "2015" (matching the textual year headers
# used on the other sheets), not a number - copy the already-textual "2015"
# label from sheet 1 and paste only its value so the type stays text while
# B1 keeps its existing header style (s=1).
$ws1.Range("B1").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").ClearFormats()
$ws6.Range("B2").Value = 630

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").ClearFormats()
$ws6.Range("B3").Value = 99

Write-Host "Headers added and data updated."
